# Commit: "Added header with case management location"
#
# In the merge-field header block:
#   <<cs_{writtenByJudge}>><<hearingLocation.site_name>><<else>> Online Civil Claims<<es_>>
# the field `hearingLocation.site_name` is renamed to `caseManagementLocation.site_name`.
# The anchor text "}>><<hearingLocation.site_name>>" is unique in the document (the other
# `hearingLocation.*` merge fields further down the template are left untouched), so a
# literal (non-wildcard) Find/Replace targets exactly this header occurrence.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$found = $find.Execute(
    "}>><<hearingLocation.site_name>>",   # FindText
    $true,                                  # MatchCase
    $false,                                 # MatchWholeWord
    $false,                                 # MatchWildcards
    $false,                                 # MatchSoundsLike
    $false,                                 # MatchAllWordForms
    $true,                                  # Forward
    1,                                       # Wrap (wdFindContinue)
    $false,                                 # Format
    "}>><<caseManagementLocation.site_name>>", # ReplaceWith
    2)                                       # Replace (wdReplaceAll)

Write-Output "Replaced hearingLocation.site_name -> caseManagementLocation.site_name: $found"
